$d = $word.ActiveDocument

# Paragraph 1 holds the bookmark-style placeholder text.
$p1 = $d.Paragraphs(1)

# Add a paragraph border (top/left/bottom/right) with 5pt spacing and no
# visible line, matching the pattern already used elsewhere in this
# document (e.g. the third paragraph).
$b = $p1.Format.Borders
$b.DistanceFromTop = 5
$b.DistanceFromLeft = 5
$b.DistanceFromBottom = 5
$b.DistanceFromRight = 5

# Widen the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1.Format.LeftIndent = 11.25

# Update the placeholder text and drop the trailing space run by replacing
# the whole "text + trailing space" span in one go.
$d.Content.Find.Execute("**ID__AFFARS_pgi_5304_topic_17__ID** ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "**ID__AFFARS_SMC_PGI_5304_803__ID**", 2)
